$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '21.752.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.539.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3877'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3194'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07202'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.056'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.24%  '
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.636'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.597'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.541.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001110'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06591'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.137'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.66%  '
$ws.Range('E23').Value = '  -5.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.388'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '21.759.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.372'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.844'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.714.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9693'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -12.89%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.915'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08184'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.906'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.90%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06071'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.138'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.485'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -16.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02204'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2037'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.191'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5746'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.08%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.746'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5518'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '117.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.865'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.30%  '
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06723'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.10%  '
